$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: only C2 changes
$ws.Range("C2").Value = 1929829825410267

# Row 3: RandomForestRegressor label unchanged, values change
$ws.Range("B3").Value = 808712703994307.4
$ws.Range("C3").Value = 831209417375311
$ws.Range("D3").Value = 1375376561930154

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 21487645573253.58
$ws.Range("C4").Value = 21660606614359.04
$ws.Range("D4").Value = 1583515623954792

# Row 5: AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 355598709141558.1
$ws.Range("C5").Value = 187770188586593.7
$ws.Range("D5").Value = 313645161278879.9
